$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 14063
$ws.Cells.Item(3, 6).Value = 13810
$ws.Cells.Item(4, 6).Value = 827
$ws.Cells.Item(8, 6).Value = 32
$ws.Cells.Item(10, 6).Value = 794
$ws.Cells.Item(11, 6).Value = 2171
$ws.Cells.Item(12, 6).Value = 159
$ws.Cells.Item(13, 6).Value = 116
$ws.Cells.Item(14, 6).Value = 94
$ws.Cells.Item(15, 6).Value = 197
$ws.Cells.Item(17, 6).Value = 573
$ws.Cells.Item(18, 6).Value = 457
$ws.Cells.Item(19, 6).Value = 487
$ws.Cells.Item(20, 6).Value = 340
$ws.Cells.Item(21, 6).Value = 24
$ws.Cells.Item(22, 6).Value = 298
$ws.Cells.Item(23, 6).Value = 869
$ws.Cells.Item(24, 6).Value = 131
$ws.Cells.Item(25, 6).Value = 59
$ws.Cells.Item(26, 6).Value = 17
$ws.Cells.Item(29, 6).Value = 57
$ws.Cells.Item(30, 6).Value = 23

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 54
$ws.Cells.Item(4, 6).Value = 142
$ws.Cells.Item(6, 6).Value = 104
$ws.Cells.Item(8, 6).Value = 1650
$ws.Cells.Item(13, 6).Value = 84
$ws.Cells.Item(15, 6).Value = 1764

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 109

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 14063
$ws.Cells.Item(4, 6).Value = 13810
$ws.Cells.Item(5, 6).Value = 827
$ws.Cells.Item(9, 6).Value = 32
$ws.Cells.Item(11, 6).Value = 794
$ws.Cells.Item(12, 6).Value = 54
$ws.Cells.Item(14, 6).Value = 2171
$ws.Cells.Item(15, 6).Value = 109
$ws.Cells.Item(16, 6).Value = 159
$ws.Cells.Item(17, 6).Value = 159
$ws.Cells.Item(18, 6).Value = 116
$ws.Cells.Item(19, 6).Value = 94
$ws.Cells.Item(20, 6).Value = 197
$ws.Cells.Item(21, 6).Value = 142
$ws.Cells.Item(24, 6).Value = 104
$ws.Cells.Item(26, 6).Value = 573
$ws.Cells.Item(27, 6).Value = 457
$ws.Cells.Item(28, 6).Value = 487
$ws.Cells.Item(29, 6).Value = 340
$ws.Cells.Item(30, 6).Value = 24
$ws.Cells.Item(31, 6).Value = 298
$ws.Cells.Item(32, 6).Value = 869
$ws.Cells.Item(34, 6).Value = 1650
$ws.Cells.Item(39, 6).Value = 131
$ws.Cells.Item(40, 6).Value = 59
$ws.Cells.Item(41, 6).Value = 17
$ws.Cells.Item(43, 6).Value = 84
$ws.Cells.Item(46, 6).Value = 57
$ws.Cells.Item(47, 6).Value = 23
$ws.Cells.Item(48, 6).Value = 1764
